$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Prime a scratch cell with Text number format so the style (xf) record
# used for forcing text-typed numeric-looking values already exists,
# and gets reused (no duplicate style entries) for every cell below.
$scratch = $ws.Range("Z200")
$scratch.NumberFormat = "@"

$ws.Range("D2").Value = "52.134.53"
$ws.Range("E2").Value = "  +0.50%  "
$ws.Range("D3").Value = "2.902.32"
$ws.Range("E3").Value = "  +3.34%  "
$ws.Range("E4").Value = "  +0.01%  "
$c = $ws.Range("D5")
$c.NumberFormat = "@"
$c.Value = "351.71"
$c.Style = "Normal"
$ws.Range("E5").Value = "  -0.54%  "
$c = $ws.Range("D6")
$c.NumberFormat = "@"
$c.Value = "113.12"
$c.Style = "Normal"
$ws.Range("E6").Value = "  +0.97%  "
$c = $ws.Range("D7")
$c.NumberFormat = "@"
$c.Value = "0.555"
$c.Style = "Normal"
$ws.Range("E7").Value = "  -0.13%  "
$ws.Range("E8").Value = "  +0.05%  "
$c = $ws.Range("D9")
$c.NumberFormat = "@"
$c.Value = "0.619"
$c.Style = "Normal"
$ws.Range("E9").Value = "  -0.89%  "
$c = $ws.Range("D10")
$c.NumberFormat = "@"
$c.Value = "39.59"
$c.Style = "Normal"
$ws.Range("E10").Value = "  -1.67%  "
$ws.Range("B11").Value = "TRON"
$ws.Range("C11").Value = "https://coinranking.com/coin/qUhEFk1I61atv+tron-trx"
$c = $ws.Range("D11")
$c.NumberFormat = "@"
$c.Value = "0.136"
$c.Style = "Normal"
$ws.Range("E11").Value = "  +0.62%  "
$ws.Range("B12").Value = "Dogecoin"
$ws.Range("C12").Value = "https://coinranking.com/coin/a91GCGd_u96cF+dogecoin-doge"
$c = $ws.Range("D12")
$c.NumberFormat = "@"
$c.Value = "0.0862"
$c.Style = "Normal"
$ws.Range("E12").Value = "  +2.66%  "
$ws.Range("E13").Value = "  -1.05%  "
$c = $ws.Range("D14")
$c.NumberFormat = "@"
$c.Value = "7.71"
$c.Style = "Normal"
$ws.Range("E14").Value = "  -0.86%  "
$ws.Range("D15").Value = "3.361.03"
$ws.Range("E15").Value = "  +3.35%  "
$ws.Range("D16").Value = "2.905.89"
$ws.Range("E16").Value = "  +3.42%  "
$ws.Range("E17").Value = "  +3.73%  "
$ws.Range("D18").Value = "52.217.94"
$ws.Range("E18").Value = "  +0.66%  "
$c = $ws.Range("D19")
$c.NumberFormat = "@"
$c.Value = "3.33"
$c.Style = "Normal"
$ws.Range("E19").Value = "  +2.82%  "
$ws.Range("E20").Value = "  -0.44%  "
$c = $ws.Range("D21")
$c.NumberFormat = "@"
$c.Value = "13.91"
$c.Style = "Normal"
$ws.Range("E21").Value = "  +2.20%  "
$ws.Range("D22").Value = "0.0₃0974"
$ws.Range("E22").Value = "  -0.18%  "
$c = $ws.Range("D23")
$c.NumberFormat = "@"
$c.Value = "71.15"
$c.Style = "Normal"
$ws.Range("E23").Value = "  +1.14%  "
$c = $ws.Range("D24")
$c.NumberFormat = "@"
$c.Value = "269.44"
$c.Style = "Normal"
$ws.Range("E24").Value = "  +0.73%  "
$ws.Range("E25").Value = "  +1.29%  "
$ws.Range("E26").Value = "  +13.52%  "
$c = $ws.Range("D27")
$c.NumberFormat = "@"
$c.Value = "26.75"
$c.Style = "Normal"
$ws.Range("E27").Value = "  +2.16%  "
$c = $ws.Range("D28")
$c.NumberFormat = "@"
$c.Value = "0.999"
$c.Style = "Normal"
$ws.Range("E28").Value = "  -0.14%  "
$c = $ws.Range("D29")
$c.NumberFormat = "@"
$c.Value = "10.59"
$c.Style = "Normal"
$ws.Range("E29").Value = "  +2.01%  "
$c = $ws.Range("D30")
$c.NumberFormat = "@"
$c.Value = "0.103"
$c.Style = "Normal"
$ws.Range("E30").Value = "  +15.35%  "
$c = $ws.Range("D31")
$c.NumberFormat = "@"
$c.Value = "6.72"
$c.Style = "Normal"
$ws.Range("E31").Value = "  +9.56%  "
$c = $ws.Range("D32")
$c.NumberFormat = "@"
$c.Value = "37.47"
$c.Style = "Normal"
$ws.Range("E32").Value = "  -3.62%  "
$c = $ws.Range("D33")
$c.NumberFormat = "@"
$c.Value = "2.27"
$c.Style = "Normal"
$ws.Range("E33").Value = "  -0.42%  "
$c = $ws.Range("D34")
$c.NumberFormat = "@"
$c.Value = "6.21"
$c.Style = "Normal"
$ws.Range("E34").Value = "  +12.17%  "
$c = $ws.Range("D35")
$c.NumberFormat = "@"
$c.Value = "53.07"
$c.Style = "Normal"
$ws.Range("E36").Value = "  -2.12%  "
$c = $ws.Range("D37")
$c.NumberFormat = "@"
$c.Value = "0.999"
$c.Style = "Normal"
$ws.Range("E37").Value = "  -0.11%  "
$ws.Range("E38").Value = "  +4.30%  "
$ws.Range("E40").Value = "  +1.28%  "
$ws.Range("E41").Value = "  +7.72%  "
$ws.Range("E42").Value = "  +1.20%  "
$c = $ws.Range("D43")
$c.NumberFormat = "@"
$c.Value = "23.03"
$c.Style = "Normal"
$ws.Range("E43").Value = "  +4.84%  "
$c = $ws.Range("D44")
$c.NumberFormat = "@"
$c.Value = "119.02"
$c.Style = "Normal"
$ws.Range("E44").Value = "  -0.35%  "
$ws.Range("E45").Value = "  -1.94%  "
$c = $ws.Range("D46")
$c.NumberFormat = "@"
$c.Value = "2.54"
$c.Style = "Normal"
$ws.Range("E46").Value = "  +3.32%  "
$ws.Range("B47").Value = "NEARProtocol"
$ws.Range("C47").Value = "https://coinranking.com/coin/DCrsaMv68+nearprotocol-near"
$c = $ws.Range("D47")
$c.NumberFormat = "@"
$c.Value = "3.50"
$c.Style = "Normal"
$ws.Range("E47").Value = "  +0.40%  "
$ws.Range("B48").Value = "Maker"
$ws.Range("C48").Value = "https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr"
$ws.Range("D48").Value = "2.171.67"
$ws.Range("E48").Value = "  +3.08%  "
$c = $ws.Range("D49")
$c.NumberFormat = "@"
$c.Value = "0.261"
$c.Style = "Normal"
$ws.Range("E49").Value = "  +20.17%  "
$c = $ws.Range("D50")
$c.NumberFormat = "@"
$c.Value = "0.0346"
$c.Style = "Normal"
$ws.Range("E50").Value = "  +9.87%  "
$c = $ws.Range("D51")
$c.NumberFormat = "@"
$c.Value = "0.953"
$c.Style = "Normal"
$ws.Range("E51").Value = "  -1.29%  "

# Clean up the scratch cell entirely so it leaves no trace in the sheet
$scratch.NumberFormat = "General"
$scratch.Style = "Normal"
$scratch.Clear()
